{"js": "// 1) Insert a new, completely empty paragraph right before the paragraph\n//    that contains \"Here is the gdrive link of pre-scraped data.\"\nconst body = context.document.body;\nconst gdriveResults = body.search(\"Here is the gdrive link of pre-scraped data.\", { matchCase: true });\ngdriveResults.load(\"items\");\nawait context.sync();\n\nif (gdriveResults.items.length > 0) {\n  const gdriveParagraph = gdriveResults.items[0].paragraphs.getFirst();\n  gdriveParagraph.insertParagraph(\"\", \"Before\");\n}\n\n// 2) Update the \"Proper Human-in-the-loop Usage\" bullet: \"(Agent should\" -> \"(E.g. Agents should\"\nconst humanLoopResults = body.search(\n  \"Proper Human-in-the-loop Usage (Agent should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\",\n  { matchCase: true }\n);\nhumanLoopResults.load(\"items\");\nawait context.sync();\n\nif (humanLoopResults.items.length > 0) {\n  humanLoopResults.items[0].insertText(\n    \"Proper Human-in-the-loop Usage (E.g. Agents should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\",\n    \"Replace\"\n  );\n}\n\n// 3) Update the \"Reasoning Capability\" bullet: add \"E.g. \" before \"This should\" and append\n//    the extra sentence about reasoning steps before the closing \"- (point 3)\".\nconst reasoningResults = body.search(\n  \"Reasoning Capability (This should be implemented as an AI reasoning agent. For development, please use non-reasoning models.) - (point 3)\",\n  { matchCase: true }\n);\nreasoningResults.load(\"items\");\nawait context.sync();\n\nif (reasoningResults.items.length > 0) {\n  reasoningResults.items[0].insertText(\n    \"Reasoning Capability (E.g. This should be implemented as an AI reasoning agent. For development, please use non-reasoning models. And reasoning steps and intermediate results should be displayed in the frontend) - (point 3)\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop script that applies the task's requested edits:\n#  1. Insert a new, completely empty paragraph right before the paragraph\n#     that contains \"Here is the gdrive link of pre-scraped data.\"\n#  2. Update the \"Proper Human-in-the-loop Usage\" bullet:\n#     \"(Agent should\" -> \"(E.g. Agents should\"\n#  3. Update the \"Reasoning Capability\" bullet:\n#     \"(This should\" -> \"(E.g. This should\" and append the extra sentence\n#     about reasoning steps/intermediate results before \"- (point 3)\".\n\n$d = $word.ActiveDocument\n\n# --- 1) Insert blank paragraph before the gdrive-link paragraph ---------\n$findRange = $d.Content\n$find = $findRange.Find\n$find.ClearFormatting()\n$find.Text = \"Here is the gdrive link of pre-scraped data.\"\n$found = $find.Execute()\nif ($found) {\n    $targetParagraph = $findRange.Paragraphs.Item(1)\n    $targetParagraph.Range.InsertParagraphBefore()\n}\n\n# --- 2) Fix \"Proper Human-in-the-loop Usage\" wording ---------------------\n$r2 = $d.Content\n$find2 = $r2.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Proper Human-in-the-loop Usage (Agent should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\"\n$find2.Replacement.Text = \"Proper Human-in-the-loop Usage (E.g. Agents should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# --- 3) Fix \"Reasoning Capability\" wording --------------------------------\n$r3 = $d.Content\n$find3 = $r3.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"Reasoning Capability (This should be implemented as an AI reasoning agent. For development, please use non-reasoning models.) - (point 3)\"\n$find3.Replacement.Text = \"Reasoning Capability (E.g. This should be implemented as an AI reasoning agent. For development, please use non-reasoning models. And reasoning steps and intermediate results should be displayed in the frontend) - (point 3)\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2) | Out-Null\n"}
